# Insert a new record row into the daily-price log for
# "Hortaliza, Terminal La Palmera de La Serena - Ajo".
#
# A new row of data (2023-11-28 / Primera / China) is inserted at row 482,
# pushing every existing row from 482-582 down by one (to 483-583). The
# workbook's used-range grows from A1:R582 to A1:R583 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 482..582 down to 483..583 by inserting a blank row at 482.
$ws.Rows(482).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(482, 1).Value  = 8
$ws.Cells.Item(482, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(482, 3).Value  = 'Coquimbo'
$ws.Cells.Item(482, 4).Value  = 45258
$ws.Cells.Item(482, 5).Value  = 4
$ws.Cells.Item(482, 6).Value  = 100112003
$ws.Cells.Item(482, 7).Value  = 'Ajo'
$ws.Cells.Item(482, 8).Value  = 'Chino'
$ws.Cells.Item(482, 9).Value  = 'Primera'
$ws.Cells.Item(482, 10).Value = 400
$ws.Cells.Item(482, 11).Value = 23000
$ws.Cells.Item(482, 12).Value = 24000
$ws.Cells.Item(482, 13).Value = 23500
$ws.Cells.Item(482, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(482, 15).Value = 'China'
$ws.Cells.Item(482, 16).Value = 2350
$ws.Cells.Item(482, 17).Value = 10
$ws.Cells.Item(482, 18).Value = 'Hortaliza'
